$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
# Match the body-text ("Normal") style used by the rest of the document
# instead of inheriting the Heading1 style from the title paragraph.
$bodyStyle = $d.Paragraphs.Item(4).Style
$metaPara.Style = $bodyStyle

$metaStart = $metaPara.Range.Start
$metaEnd = $metaPara.Range.End
$metaTextRange = $d.Range($metaStart, $metaEnd - 1)
$metaTextRange.Text = "Meta description: Discover the Dragon Emperor online slot game from Aristocrat, play for free, win with numerous betting options & free spins up to 15 and 40X multiplier."

# Bold just the "Meta description" label (first 17 characters).
$boldRange = $d.Range($metaStart, $metaStart + 16)
$boldRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicated bold "Play the Free Dragon Emperor Slot by
#    Aristocrat - Review" paragraph near the end of the document, and
# 3) Replace the text of the italic paragraph that followed it with the
#    new image-prompt text (keeping its italic formatting).
# ------------------------------------------------------------------
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Bold -and ($p.Range.Text.TrimEnd() -eq "Play the Free Dragon Emperor Slot by Aristocrat - Review")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}

$italicPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd().StartsWith("Discover the Dragon Emperor online slot game from Aristocrat")) {
        $italicPara = $p
        break
    }
}

$newPromptText = "Create an image to capture the essence of the Dragon Emperor online slot game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing in front of a majestic dragon, holding a treasure chest with a big smile on their face. The backdrop should showcase the magical and mysterious atmosphere of the game. Use bright and vibrant colors to make the image pop and catch players' attention. The image should showcase the thrill and excitement of the game, encouraging players to embark on the journey to uncover the Dragon's treasure."

if ($italicPara -ne $null) {
    $s = $italicPara.Range.Start
    $e = $italicPara.Range.End
    $textRange = $d.Range($s, $e - 1)
    $textRange.Text = $newPromptText
}

Write-Host "Edit complete"
